# Delete row 3 (Manal ElMetwally ElHoseiny) from the Employees sheet.
# Excel shifts the remaining rows up automatically, matching the target
# workbook where row 3 is removed and rows 4-28 become rows 3-27.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(3).Delete()
